# Weekly data refresh: insert a new observation as row 149, pushing the
# existing rows (149-220) down by one (220 -> 221).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 149 (shifts 149..220 down to 150..221).
$ws.Rows(149).Insert()

# Populate the newly inserted row 149 with the new weekly record.
$ws.Range("A149").Value = 8
$ws.Range("B149").Value = "Terminal La Palmera de La Serena"
$ws.Range("C149").Value = "Coquimbo"
$ws.Range("D149").Value = 45016
$ws.Range("E149").Value = 4
$ws.Range("F149").Value = 100112001
$ws.Range("G149").Value = "Berenjena"
$ws.Range("H149").Value = "Sin especificar"
$ws.Range("I149").Value = "Primera"
$ws.Range("J149").Value = 360
$ws.Range("K149").Value = 11000
$ws.Range("L149").Value = 12000
$ws.Range("M149").Value = 11500
$ws.Range("N149").Value = "$/caja 50 unidades"
$ws.Range("O149").Value = "Región de Arica y Parinacota"
$ws.Range("P149").Value = 230
$ws.Range("Q149").Value = 50
$ws.Range("R149").Value = "Hortaliza"
